$d = $word.ActiveDocument

# Change 1: Front-end Technologies line
#   "React, Redux, jQuery" -> "React, " + "React Native, " + "Redux, jQuery"
$f1 = $d.Content.Find.Execute("React, Redux, jQuery", $false, $false, $false, $false, $false, $true, 1, $false, "React, React Native, Redux, jQuery", 2)

# Change 2: Back-end Technologies line
#   " SQL, RESTful APIs" -> " SQL, Apollo, RESTful APIs"
$f2 = $d.Content.Find.Execute(" SQL, RESTful APIs", $false, $false, $false, $false, $false, $true, 1, $false, " SQL, Apollo, RESTful APIs", 2)
#   ", Python, Django." -> ", Python, Django, Axios."
$f3 = $d.Content.Find.Execute(", Python, Django.", $false, $false, $false, $false, $false, $true, 1, $false, ", Python, Django, Axios.", 2)

# Change 3: Deployment and Hosting line
#   " AWS, Heroku" -> " Expo, EAS, AWS, Heroku"
$f4 = $d.Content.Find.Execute(" AWS, Heroku", $false, $false, $false, $false, $false, $true, 1, $false, " Expo, EAS, AWS, Heroku", 2)

# Change 4: Tools line
#   ", Docker" -> ", Docker."
$f5 = $d.Content.Find.Execute(", Docker", $false, $false, $false, $false, $false, $true, 1, $false, ", Docker.", 2)

Write-Output "f1=$f1 f2=$f2 f3=$f3 f4=$f4 f5=$f5"
